# "Diamons" -> "Mazas y Catapultas": the currently-playing track name in the
# sample Data sheet was corrected once the app actually hit "play". Re-fit the
# column widths to the sheet's content (as Excel does on a manual
# AutoFit-columns pass) and land the selection on D3, ready for the next row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

$ws.Range("D2").Value = "Mazas y Catapultas"

$ws.Range("D3").Select() | Out-Null
